# Psykologsystemer workbook update (02-12-2025 -> 05-12-2025).
#
# The sheet title carries the "last updated" date for this approved-systems
# list. Renaming the worksheet also re-points the `Psykologsystemer` defined
# name at the renamed sheet, which is exactly what the authoritative diff
# shows (sheet name + definedName both flip from "02-12-2025" to
# "05-12-2025", same $A$1:$G$14 extent).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newName = "Opdateret d. 05-12-2025"

if (-not ($ws.Name -eq $newName)) {
    $ws.Name = $newName
}
